# Generate Report for Handoff
# Inserts a new "756e0514-3dd6-4c71-8197-12542bec0e46" handoff entry as row 2
# on every data sheet (Overview, zh-cn, de-de), pushing the existing
# "cdf7d01a-772b-45bb-b21b-218f2a4edb3a" entry down to row 3.

$wb = $excel.ActiveWorkbook

$newGuid = "756e0514-3dd6-4c71-8197-12542bec0e46"
$oldGuid = "cdf7d01a-772b-45bb-b21b-218f2a4edb3a"

$newMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/5b14e624a0edf4adad337e0f7adcee966c79d55d/e2e/$newGuid.md"
$oldMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/5b14e624a0edf4adad337e0f7adcee966c79d55d/e2e/$oldGuid.md"

$newZhUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c9b104d844130d222eb0de64148b63dd225cff71/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newGuid.aa50abacea9f47958792410feec0f185133fb1df.zh-cn.xlf"
$oldZhUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c9b104d844130d222eb0de64148b63dd225cff71/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.9c25ab04bb3ab0acc75c841c7a2a6c42bde49f2a.zh-cn.xlf"

$newDeUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef0fcc737b1e548d1ca2940f39aecb817e7cbd32/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newGuid.aa50abacea9f47958792410feec0f185133fb1df.de-de.xlf"
$oldDeUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef0fcc737b1e548d1ca2940f39aecb817e7cbd32/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.9c25ab04bb3ab0acc75c841c7a2a6c42bde49f2a.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Push the existing data row (row 2) down to row 3, carrying its values,
# number formats and hyperlinked-cell styling with it.
$ws1.Rows.Item(2).Insert()

# Rebuild hyperlinks: the old ones still point at the pre-insert locations.
$ws1.Hyperlinks.Delete()

$ws1.Range("A2").Value2 = "$newGuid.md"
$ws1.Range("B2").Value2 = "Ready for handoff"
$ws1.Range("C2").Value2 = "Ready for handoff"
$ws1.Range("D2").Value2 = "2016-26-12 08:26:36"
$ws1.Range("A2").Style = "HyperLink"

$ws1.Hyperlinks.Add($ws1.Range("A2"), $newMdUrl, "", "", "$newGuid.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), $oldMdUrl, "", "", "$oldGuid.md")

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(2).Insert()
$ws2.Hyperlinks.Delete()

$ws2.Range("A2").Value2 = "$newGuid.md"
$ws2.Range("B2").Value2 = ".md"
$ws2.Range("C2").Value2 = "Ready for handoff"
$ws2.Range("D2").Value2 = "$newGuid.aa50abacea9f47958792410feec0f185133fb1df.zh-cn.xlf"
$ws2.Range("E2").Value2 = "2016-03-12 08:26:31"
$ws2.Range("H2").Value2 = "0001-01-01 00:00:00"
$ws2.Range("I2").Value2 = "Include"

$ws2.Range("A2").Style = "HyperLink"
$ws2.Range("B2").Style = "HyperLink"
$ws2.Range("D2").Style = "HyperLink"

$ws2.Hyperlinks.Add($ws2.Range("A2"), $newMdUrl, "", "", "$newGuid.md")
$ws2.Hyperlinks.Add($ws2.Range("B2"), $newMdUrl, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), $newZhUrl, "", "", "$newGuid.aa50abacea9f47958792410feec0f185133fb1df.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A3"), $oldMdUrl, "", "", "$oldGuid.md")
$ws2.Hyperlinks.Add($ws2.Range("B3"), $oldMdUrl, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), $oldZhUrl, "", "", "$oldGuid.9c25ab04bb3ab0acc75c841c7a2a6c42bde49f2a.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(2).Insert()
$ws3.Hyperlinks.Delete()

$ws3.Range("A2").Value2 = "$newGuid.md"
$ws3.Range("B2").Value2 = ".md"
$ws3.Range("C2").Value2 = "Ready for handoff"
$ws3.Range("D2").Value2 = "$newGuid.aa50abacea9f47958792410feec0f185133fb1df.de-de.xlf"
$ws3.Range("E2").Value2 = "2016-03-12 08:26:36"
$ws3.Range("H2").Value2 = "0001-01-01 00:00:00"
$ws3.Range("I2").Value2 = "Include"

$ws3.Range("A2").Style = "HyperLink"
$ws3.Range("B2").Style = "HyperLink"
$ws3.Range("D2").Style = "HyperLink"

$ws3.Hyperlinks.Add($ws3.Range("A2"), $newMdUrl, "", "", "$newGuid.md")
$ws3.Hyperlinks.Add($ws3.Range("B2"), $newMdUrl, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), $newDeUrl, "", "", "$newGuid.aa50abacea9f47958792410feec0f185133fb1df.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A3"), $oldMdUrl, "", "", "$oldGuid.md")
$ws3.Hyperlinks.Add($ws3.Range("B3"), $oldMdUrl, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), $oldDeUrl, "", "", "$oldGuid.9c25ab04bb3ab0acc75c841c7a2a6c42bde49f2a.de-de.xlf")
